# Apply the "fix issues de compilation" update to the Metadata sheet:
#  - Update the Date property value
#  - Insert a new "Jurisdiction" property row (with an empty value) right
#    after the "Contact" row, pushing the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Date" property value (row 8, column B).
$ws.Range("B8").Value = "2024-09-12T14:01:50+00:00"

# 2) Insert a new row above the current row 12 ("Purpose"), which becomes
#    the new row 11, directly below "Contact" (row 10).
$ws.Rows.Item(11).Insert()

# Copy the formatting from the row below (now row 12, "Description") onto
# the newly inserted row so it keeps the same style as the rest of the
# table instead of picking up a blank default style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Fill in the new "Jurisdiction" property with an empty value.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
